$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, pushing existing rows 27-91 down to 28-92
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly price report record
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(27, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value = "Metropolitana"
$ws.Cells.Item(27, 4).Value = 44544
$ws.Cells.Item(27, 5).Value = 13
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = "Arveja Verde"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 52
$ws.Cells.Item(27, 11).Value = 15000
$ws.Cells.Item(27, 12).Value = 17000
$ws.Cells.Item(27, 13).Value = 16000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 640
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
